# Staging.MilestoneType.xlsx regen: a new "BusinessKey" column was added as
# the first data column on the header row, pushing the existing Code /
# MilestoneTypeID / Name headers one column to the right (A2:D2 becomes
# BusinessKey, Code, MilestoneTypeID, Name). Row 1 (the "For internal use
# only" banner) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing header values one column to the right (write from the
# rightmost cell backwards so we don't clobber a value before reading it),
# then drop the new column's label into A2. All three existing cells keep
# their "s=1" (bold+underline) custom formatting automatically since D2/C2/B2
# get set the same way the template's other header cells already are.
# (Use .Value2 to read — .Value's getter doesn't resolve to the scalar here.)
$ws.Range("D2").Value = $ws.Range("C2").Value2
$ws.Range("C2").Value = $ws.Range("B2").Value2
$ws.Range("B2").Value = $ws.Range("A2").Value2
$ws.Range("A2").Value = "BusinessKey"
